# Append two new sheets ("area_lores_basic" and "area_pop_sum_basic") at the
# end of the workbook, holding the analogous "basic" (reduced) dataset
# statistics. The pre-existing "area_lores" / "area_pop_sum" sheets are left
# untouched.
#
# Each new sheet is created by duplicating its corresponding existing sheet
# (so it inherits identical page setup / formatting / header styling) and
# then overwriting only the numeric "value" column with the refreshed
# "_basic" statistics. The text labels in column A are left exactly as
# duplicated, which also avoids Excel re-interpreting strings such as "25%"
# as percentage numbers.

$wb = $excel.ActiveWorkbook

$wsLores = $wb.Worksheets.Item("area_lores")
$wsPopSum = $wb.Worksheets.Item("area_pop_sum")

# --- add area_lores_basic sheet (duplicate of area_lores, placed at the end) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLores.Copy($null, $lastSheet)
$wsLoresBasic = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLoresBasic.Name = "area_lores_basic"

$wsLoresBasic.Range("B2").Value = 45
$wsLoresBasic.Range("B3").Value = 5.61782385972893
$wsLoresBasic.Range("B4").Value = 5.99425330634962
$wsLoresBasic.Range("B5").Value = 0.001830784076881646
$wsLoresBasic.Range("B6").Value = 1.350757262696783
$wsLoresBasic.Range("B7").Value = 3.571077202539338
$wsLoresBasic.Range("B8").Value = 7.698480299250822
$wsLoresBasic.Range("B9").Value = 25.07158166439542

# --- add area_pop_sum_basic sheet (duplicate of area_pop_sum, placed at the end) ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPopSum.Copy($null, $lastSheet2)
$wsPopSumBasic = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPopSumBasic.Name = "area_pop_sum_basic"

$wsPopSumBasic.Range("B2").Value = 252.8020736878019
$wsPopSumBasic.Range("B3").Value = 704016
$wsPopSumBasic.Range("B4").Value = 2784.85057392933

# restore original active sheet/selection state
$wsLores.Activate()
